# Fruta / hortaliza, semanal
#
# A new weekly observation is inserted as row 63 (pushing the existing
# rows 63-145 down to 64-146). This mirrors the author's edit: inserting
# a single new price record above the row that used to be "D63".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 63; everything below shifts down by one.
$ws.Rows("63:63").Insert()

# Populate the newly inserted row with the new observation.
$ws.Range("A63").Value = 5
$ws.Range("B63").Value = "Macroferia Regional de Talca"
$ws.Range("C63").Value = "Maule"
$ws.Range("D63").Value = 45174
$ws.Range("E63").Value = 7
$ws.Range("F63").Value = 100112013
$ws.Range("G63").Value = "Alcachofa"
$ws.Range("H63").Value = "Madrigal"
$ws.Range("I63").Value = "Primera"
$ws.Range("J63").Value = 300
$ws.Range("K63").Value = 10000
$ws.Range("L63").Value = 10000
$ws.Range("M63").Value = 10000
$ws.Range("N63").Value = "$/caja 40 unidades"
$ws.Range("O63").Value = "Provincia del Elquí"
$ws.Range("P63").Value = 250
$ws.Range("Q63").Value = 40
$ws.Range("R63").Value = "Hortaliza"
